$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "732a2f3f195b8d0903914c644b645b13"
$ws.Range("B2").Value = "xuzm02"
$ws.Range("A3").Value = "b46760af1959e5a87969966dc65b6c8f"
$ws.Range("B3").Value = "xrad04"
$ws.Range("A4").Value = "a8c92e10b44148d74bc2ee990960efb2"
$ws.Range("B4").Value = "xoha05"
$ws.Range("A5").Value = "1c2ab786ac4650ebe381ddb193297cc1"
$ws.Range("B5").Value = "xtea00"
$ws.Range("A6").Value = "aff5acae4320c18a387987dedf6a0161"
$ws.Range("B6").Value = "xoja00"
$ws.Range("A7").Value = "dce52e05f6c505af6829a223db4fc7da"
$ws.Range("B7").Value = "xanm32"
$ws.Range("A8").Value = "2a6841403b63a8907284469b10fd2f1a"
$ws.Range("B8").Value = "xulj07"
$ws.Range("A9").Value = "b9fbc883c4d2a5a45c2e1e3ef35a6b1d"
$ws.Range("B9").Value = "xand03"
$ws.Range("A10").Value = "d3b875048a7eb4ac828a5a6dd78b1c62"
$ws.Range("B10").Value = "xeda07"
$ws.Range("A11").Value = "efb57303ac40a4b5d636208977032627"
$ws.Range("B11").Value = "xorv03"
$ws.Range("A12").Value = "72a7371fd7b592eb988fb6b3f8547841"
$ws.Range("B12").Value = "xubd06"
$ws.Range("A13").Value = "14687139b1b4634cd81e240500db5546"
$ws.Range("B13").Value = "xejm04"
$ws.Range("A14").Value = "ce5cc509c3022c379afe0886f6100312"
$ws.Range("B14").Value = "xvep00"
$ws.Range("A15").Value = "740bfee7b07d213852ac7ffbba842625"
$ws.Range("B15").Value = "xcht09"
$ws.Range("A16").Value = "57c7f0a500f4f1e66bc949d6a6688828"
$ws.Range("B16").Value = "xopl11"
$ws.Range("A17").Value = "14eebca56a282a3140386ce10f02c37c"
$ws.Range("B17").Value = "xtoj00"
$ws.Range("A18").Value = "b2d7f3ccf127550448da997562e96207"
$ws.Range("B18").Value = "xopa00"
$ws.Range("A19").Value = "8c09e463211cd22639e7044c6ff55713"
$ws.Range("B19").Value = "xtaj19"
$ws.Range("A20").Value = "8234ded31f4854b527628a356b6941ec"
$ws.Range("B20").Value = "xabv03"
$ws.Range("A21").Value = "2ff32b7998e9722925b7ed7bac3254a6"
$ws.Range("B21").Value = "xtrr07"
$ws.Range("A22").Value = "b76fca21389ed13baae86cd5c334cc41"
$ws.Range("B22").Value = "xloa00"
$ws.Range("A23").Value = "0d20e2813565e040a3c6a65757ef1ddb"
$ws.Range("B23").Value = "xasj06"
$ws.Range("A24").Value = "33abca8ee042c7a3fbf8889150cffd90"
$ws.Range("B24").Value = "xamj01"
$ws.Range("A25").Value = "d993140c0a2316fd3ed617b9c0daa4e2"
$ws.Range("B25").Value = "xoxm00"
$ws.Range("A26").Value = "f1e03417a159177efeb664e0331835ae"
$ws.Range("B26").Value = "xodm10"
$ws.Range("A27").Value = "5573e5b42f52ede65978037777af6b00"
$ws.Range("B27").Value = "xuzj11"
$ws.Range("A28").Value = "dce52e05f6c505af6829a223db4fc7da"
$ws.Range("B28").Value = "xanm32"
$ws.Range("A29").Value = "ce5cc509c3022c379afe0886f6100312"
$ws.Range("B29").Value = "xvep00"
$ws.Range("A30").Value = "1c2ab786ac4650ebe381ddb193297cc1"
$ws.Range("B30").Value = "xtea00"
$ws.Range("A31").Value = "72a7371fd7b592eb988fb6b3f8547841"
$ws.Range("B31").Value = "xubd06"
$ws.Range("A32").Value = "efb57303ac40a4b5d636208977032627"
$ws.Range("B32").Value = "xorv03"
$ws.Range("A33").Value = "14687139b1b4634cd81e240500db5546"
$ws.Range("B33").Value = "xejm04"
$ws.Range("A34").Value = "a8c92e10b44148d74bc2ee990960efb2"
$ws.Range("B34").Value = "xoha05"
$ws.Range("A35").Value = "d3b875048a7eb4ac828a5a6dd78b1c62"
$ws.Range("B35").Value = "xeda07"
$ws.Range("A36").Value = "2a6841403b63a8907284469b10fd2f1a"
$ws.Range("B36").Value = "xulj07"
$ws.Range("A37").Value = "b9fbc883c4d2a5a45c2e1e3ef35a6b1d"
$ws.Range("B37").Value = "xand03"
$ws.Range("A38").Value = "f83b3b203ca661fa109ff4869d0436f5"
$ws.Range("B38").Value = "xedj00"
$ws.Range("A39").Value = "b46760af1959e5a87969966dc65b6c8f"
$ws.Range("B39").Value = "xrad04"
$ws.Range("A40").Value = "57c7f0a500f4f1e66bc949d6a6688828"
$ws.Range("B40").Value = "xopl11"
$ws.Range("A41").Value = "732a2f3f195b8d0903914c644b645b13"
$ws.Range("B41").Value = "xuzm02"
$ws.Range("A42").Value = "aff5acae4320c18a387987dedf6a0161"
$ws.Range("B42").Value = "xoja00"
$ws.Range("A43").Value = "740bfee7b07d213852ac7ffbba842625"
$ws.Range("B43").Value = "xcht09"
$ws.Range("A44").Value = "d993140c0a2316fd3ed617b9c0daa4e2"
$ws.Range("B44").Value = "xoxm00"
$ws.Range("A45").Value = "b2d7f3ccf127550448da997562e96207"
$ws.Range("B45").Value = "xopa00"
$ws.Range("A46").Value = "b76fca21389ed13baae86cd5c334cc41"
$ws.Range("B46").Value = "xloa00"
$ws.Range("A47").Value = "14eebca56a282a3140386ce10f02c37c"
$ws.Range("B47").Value = "xtoj00"
$ws.Range("A48").Value = "33abca8ee042c7a3fbf8889150cffd90"
$ws.Range("B48").Value = "xamj01"
$ws.Range("A49").Value = "2ff32b7998e9722925b7ed7bac3254a6"
$ws.Range("B49").Value = "xtrr07"
$ws.Range("A50").Value = "8c09e463211cd22639e7044c6ff55713"
$ws.Range("B50").Value = "xtaj19"
$ws.Range("A51").Value = "0d20e2813565e040a3c6a65757ef1ddb"
$ws.Range("B51").Value = "xasj06"
$ws.Range("A52").Value = "f1e03417a159177efeb664e0331835ae"
$ws.Range("B52").Value = "xodm10"
$ws.Range("A53").Value = "5573e5b42f52ede65978037777af6b00"
$ws.Range("B53").Value = "xuzj11"
$ws.Range("A54").Value = "8234ded31f4854b527628a356b6941ec"
$ws.Range("B54").Value = "xabv03"
$ws.Range("A55").Value = "2a6841403b63a8907284469b10fd2f1a"
$ws.Range("B55").Value = "xulj07"
$ws.Range("A56").Value = "efb57303ac40a4b5d636208977032627"
$ws.Range("B56").Value = "xorv03"
$ws.Range("A57").Value = "732a2f3f195b8d0903914c644b645b13"
$ws.Range("B57").Value = "xuzm02"
$ws.Range("A58").Value = "ce5cc509c3022c379afe0886f6100312"
$ws.Range("B58").Value = "xvep00"
$ws.Range("A59").Value = "72a7371fd7b592eb988fb6b3f8547841"
$ws.Range("B59").Value = "xubd06"
$ws.Range("A60").Value = "740bfee7b07d213852ac7ffbba842625"
$ws.Range("B60").Value = "xcht09"
$ws.Range("A61").Value = "14687139b1b4634cd81e240500db5546"
$ws.Range("B61").Value = "xejm04"
$ws.Range("A62").Value = "f83b3b203ca661fa109ff4869d0436f5"
$ws.Range("B62").Value = "xedj00"
$ws.Range("A63").Value = "aff5acae4320c18a387987dedf6a0161"
$ws.Range("B63").Value = "xoja00"
$ws.Range("A64").Value = "b9fbc883c4d2a5a45c2e1e3ef35a6b1d"
$ws.Range("B64").Value = "xand03"
$ws.Range("A65").Value = "d3b875048a7eb4ac828a5a6dd78b1c62"
$ws.Range("B65").Value = "xeda07"
$ws.Range("A66").Value = "57c7f0a500f4f1e66bc949d6a6688828"
$ws.Range("B66").Value = "xopl11"
$ws.Range("A67").Value = "1c2ab786ac4650ebe381ddb193297cc1"
$ws.Range("B67").Value = "xtea00"
$ws.Range("A68").Value = "b46760af1959e5a87969966dc65b6c8f"
$ws.Range("B68").Value = "xrad04"
$ws.Range("A69").Value = "dce52e05f6c505af6829a223db4fc7da"
$ws.Range("B69").Value = "xanm32"
$ws.Range("A70").Value = "8234ded31f4854b527628a356b6941ec"
$ws.Range("B70").Value = "xabv03"
$ws.Range("A71").Value = "b76fca21389ed13baae86cd5c334cc41"
$ws.Range("B71").Value = "xloa00"
$ws.Range("A72").Value = "33abca8ee042c7a3fbf8889150cffd90"
$ws.Range("B72").Value = "xamj01"
$ws.Range("A73").Value = "2ff32b7998e9722925b7ed7bac3254a6"
$ws.Range("B73").Value = "xtrr07"
$ws.Range("A74").Value = "b2d7f3ccf127550448da997562e96207"
$ws.Range("B74").Value = "xopa00"
$ws.Range("A75").Value = "14eebca56a282a3140386ce10f02c37c"
$ws.Range("B75").Value = "xtoj00"
$ws.Range("A76").Value = "0d20e2813565e040a3c6a65757ef1ddb"
$ws.Range("B76").Value = "xasj06"
$ws.Range("A77").Value = "5573e5b42f52ede65978037777af6b00"
$ws.Range("B77").Value = "xuzj11"
$ws.Range("A78").Value = "8c09e463211cd22639e7044c6ff55713"
$ws.Range("B78").Value = "xtaj19"
